# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country names (ranking shuffled as totals were refreshed) ---

# Niger / Georgia swap places (rows 143-144)
$ws.Range("A143").Value = "Georgia"
$ws.Range("A144").Value = "Niger"

# Guadalupe / Islas Caimanes / Islas Feroe rotate (rows 176-178):
# Islas Feroe moves to the top of the trio
$ws.Range("A176").Value = "Islas Feroe"
$ws.Range("A177").Value = "Guadalupe"
$ws.Range("A178").Value = "Islas Caimanes"

# Islas Malvinas / Groenlandia swap places (rows 210-211)
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Update the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 09:23"

# --- Update refreshed case statistics ---

# row 6
$ws.Range("B6").Value = 1389097
$ws.Range("C6").Value = 3603
$ws.Range("D6").Value = 887295
$ws.Range("E6").Value = 469675
$ws.Range("G6").Value = 31
$ws.Range("H6").Value = 32127

# row 51
$ws.Range("E51").Value = 3404
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 138

# row 56
$ws.Range("D56").Value = 20388
$ws.Range("E56").Value = 11148
$ws.Range("G56").Value = 28
$ws.Range("H56").Value = 1277

# row 101
$ws.Range("B101").Value = 4435
$ws.Range("C101").Value = 11
$ws.Range("D101").Value = 3329
$ws.Range("E101").Value = 510

# row 143 (now Georgia) - new figures
$ws.Range("B143").Value = 1131
$ws.Range("C143").Value = 14
$ws.Range("D143").Value = 920
$ws.Range("E143").Value = 195
$ws.Range("H143").Value = 16

# row 144 (now Niger) - inherits the previous Niger row's figures
$ws.Range("B144").Value = 1124
$ws.Range("D144").Value = 1025
$ws.Range("E144").Value = 30
$ws.Range("H144").Value = 69

# row 176 (now Islas Feroe) - new figures
$ws.Range("B176").Value = 214
$ws.Range("C176").Value = 22
$ws.Range("D176").Value = 188
$ws.Range("E176").Value = 26
$ws.Range("H176").Value = 0

# row 177 (now Guadalupe) - inherits the previous Guadalupe row's figures
$ws.Range("D177").Value = 176
$ws.Range("E177").Value = 13
$ws.Range("H177").Value = 14

# row 178 (now Islas Caimanes) - inherits the previous Islas Caimanes row's figures
$ws.Range("B178").Value = 203
$ws.Range("D178").Value = 202
$ws.Range("E178").Value = 0
$ws.Range("H178").Value = 1
